# Removed Test Case Inter-Dependency
# - Update product name to a unique "-1st" suffixed value (on both sheets, since
#   they both display the same product name text)
# - Change the "shortname" value on ProductLoanInput from numeric 4351 to text "435a"
# - Make ProductLoanInput the active/selected sheet instead of ProductLoanOutput
# - Tighten the selection on ProductLoanInput from B2:B3 down to just B3

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name on both sheets (they both hold the same text)
$wsInput.Range("B1").Value = "4351-Simple-Group-Loan-Product-Loanproduct-OVERDUEFEEFLAT-1st"
$wsOutput.Range("B1").Value = "4351-Simple-Group-Loan-Product-Loanproduct-OVERDUEFEEFLAT-1st"

# Change shortname cell to a text value
$wsInput.Range("B2").Value = "435a"

# Make ProductLoanInput the active sheet (was ProductLoanOutput)
$wsInput.Activate()

# Adjust selection on input sheet
$wsInput.Range("B3").Select()
